$d = $word.ActiveDocument

$replacements = @(
    @("636×4=2544", "237×3=711"),
    @("901×5=4505", "502×9=4518"),
    @("245×7=1715", "147×6=882"),
    @("206×6=1236", "311×2=622"),
    @("771×2=1542", "646×7=4522"),
    @("172×2=344", "888×8=7104"),
    @("144×8=1152", "177×7=1239"),
    @("880×6=5280", "276×4=1104"),
    @("527×3=1581", "126×3=378"),
    @("725×8=5800", "494×5=2470"),
    @("707×9=6363", "365×7=2555"),
    @("914×7=6398", "636×8=5088"),
    @("729×4=2916", "138×8=1104"),
    @("176×3=528", "772×5=3860"),
    @("220×5=1100", "308×8=2464"),
    @("423×8=3384", "905×2=1810"),
    @("151×6=906", "367×6=2202"),
    @("746×9=6714", "246×7=1722"),
    @("632×5=3160", "413×3=1239"),
    @("680×2=1360", "875×7=6125"),
    @("666×3=1998", "596×9=5364"),
    @("402×3=1206", "693×5=3465"),
    @("392×6=2352", "311×9=2799"),
    @("645×5=3225", "938×9=8442"),
    @("197×8=1576", "618×2=1236")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
